# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados" timestamp text in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 09:50"

# Israel (row 21) - updated counts
$ws.Range("B21").Value = 7030
$ws.Range("C21").Value = 173
$ws.Range("D21").Value = 338
$ws.Range("E21").Value = 6656
$ws.Range("F21").Value = 115
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 36

# Rusia moves up in the ranking (now row 25) with fresh data
$ws.Range("A25").Value = "Rusia"
$ws.Range("B25").Value = 4149
$ws.Range("C25").Value = 601
$ws.Range("D25").Value = 281
$ws.Range("E25").Value = 3834
$ws.Range("F25").Value = 8
$ws.Range("G25").Value = 4
$ws.Range("H25").Value = 34

# Chequia shifts down to row 26 (data unchanged)
$ws.Range("A26").Value = "Chequia"
$ws.Range("B26").Value = 3869
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 71
$ws.Range("E26").Value = 3752
$ws.Range("F26").Value = 77
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 46

# Irlanda shifts down to row 27 (data unchanged)
$ws.Range("A27").Value = "Irlanda"
$ws.Range("B27").Value = 3849
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 5
$ws.Range("E27").Value = 3746
$ws.Range("F27").Value = 109
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 98

# Bosnia y Herzegovina (row 71) - updated counts
$ws.Range("D71").Value = 27
$ws.Range("E71").Value = 500
